$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B9").NumberFormat = "@"

$ws.Range("B10").NumberFormat = "@"
$ws.Range("A10").Value = "Zero"
$ws.Range("B10").Value = "0002134"
